$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (ID = 6): set Usuario (B7) to "123" (forced as text, like the admin
# typing the table's reservation/phone code) and Estado (C7) to "Ocupada"
$ws.Range("B7").Value = "'123"
$ws.Range("C7").Value = "Ocupada"

# Row 8 (ID = 8): clear the empty Usuario (B8) cell so it is no longer present
$ws.Range("B8").ClearContents()
